$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: add birth year for Dương Văn Miếu's wife entry ---
$ws.Range("E2").Value = 1910

# --- Row 3: add birth year, clear bold formatting (was inheriting bold style) ---
$ws.Range("E3").Value = 1943
$ws.Range("E3").NumberFormat = "General"

# --- Row 3: fill in father ID, make it bold+centered like other ID cells ---
$ws.Range("H3").Value = 140001
$ws.Range("H3").Font.Bold = $true
$ws.Range("H3").HorizontalAlignment = -4108

# --- Rows 6 & 7: the two people were swapped (data moved between rows) ---
# Row 6 becomes "140001g2" / Dương Thị Chinh, with a new birth year and no "x" marker
$ws.Range("A6").Value = "140001g2"
$ws.Range("C6").Clear()
$ws.Range("D6").Value = "Dương Thị Chinh"
$ws.Range("E6").Value = 1950

# Row 7 becomes 150002 / Dương Minh Tự, keeps the "x" marker, adds birth year + father ID
$ws.Range("A7").Value = 150002
$ws.Range("C7").Value = "x"
$ws.Range("D7").Value = "Dương Minh Tự"
$ws.Range("E7").Value = 1953
$ws.Range("H7").Value = 140001

# --- Selection moved to E3 ---
$ws.Range("E3").Select()
